$wb = $excel.ActiveWorkbook

# --- Transactions sheet: it is no longer the active/selected tab ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$null = $wsTrans.Range("F19").Select()

# --- Repayment Schedule sheet edits ---
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column before O (column 15), shifting old O/P -> P/Q
$ws.Columns.Item(15).Insert()

# The old "Late" header (was in N1) moves to the newly inserted O1 column;
# N1 becomes blank.
$ws.Range("O1").Value = $ws.Range("N1").Value()
$ws.Range("N1").ClearContents()

# Clear the N column data values (rows 3-14) - they become blank cells,
# and populate the new O column with 0 for those rows.
for ($r = 3; $r -le 14; $r++) {
    $ws.Cells.Item($r, 14).ClearContents()
    $ws.Cells.Item($r, 15).Value = 0
}

# Make Repayment Schedule the active (selected) tab, with the new selection.
$ws.Activate()
$null = $ws.Range("I19").Select()
